# Auto-generated edit script applying the Louisoix_Profits profit-column refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled price-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3514.1428
$ws.Range("J43").Value = 3720
$ws.Range("L43").Value = 3720
$ws.Range("N43").Value = -3858
$ws.Range("H74").Value = 7374.25
$ws.Range("I74").Value = 6499
$ws.Range("K74").Value = 6499
$ws.Range("M74").Value = -5563
$ws.Range("H77").Value = 7374.25
$ws.Range("I77").Value = 6499
$ws.Range("K77").Value = 32495
$ws.Range("M77").Value = -27815
$ws.Range("H88").Value = 1995.1666
$ws.Range("J88").Value = 1476.1428
$ws.Range("L88").Value = 1476.1428
$ws.Range("N88").Value = -2288.1428
$ws.Range("H91").Value = 1995.1666
$ws.Range("J91").Value = 1476.1428
$ws.Range("L91").Value = 1476.1428
$ws.Range("N91").Value = -4284.1428
$ws.Range("H103").Value = 601.25
$ws.Range("J103").Value = 742
$ws.Range("L103").Value = 2226
$ws.Range("N103").Value = -3398
$ws.Range("H106").Value = 16962.666
$ws.Range("I106").Value = 10444
$ws.Range("K106").Value = 10444
$ws.Range("M106").Value = -9813
$ws.Range("H135").Value = 1654.0834
$ws.Range("I135").Value = 1049.4706
$ws.Range("K135").Value = 9445.235400000001
$ws.Range("M135").Value = -6910.235400000001
$ws.Range("H137").Value = 6850.4546
$ws.Range("I137").Value = 1365.5
$ws.Range("J137").Value = 9984.714
$ws.Range("K137").Value = 4096.5
$ws.Range("L137").Value = 29954.142
$ws.Range("M137").Value = -1546.5
$ws.Range("N137").Value = -35054.142
$ws.Range("H141").Value = 6556.3184
$ws.Range("I141").Value = 6911.95
$ws.Range("K141").Value = 20735.85
$ws.Range("M141").Value = -15555.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1450
$ws.Range("J4").Value = 1450
$ws.Range("L4").Value = 1450
$ws.Range("N4").Value = -1682
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H32").Value = 27992.574
$ws.Range("I32").Value = 29833
$ws.Range("K32").Value = 29833
$ws.Range("M32").Value = -29546
$ws.Range("H61").Value = 3996.182
$ws.Range("I61").Value = 3691.8
$ws.Range("K61").Value = 3691.8
$ws.Range("M61").Value = -3479.8
$ws.Range("H110").Value = 3098.963
$ws.Range("I110").Value = 3325.7896
$ws.Range("J110").Value = 2560.25
$ws.Range("K110").Value = 3325.7896
$ws.Range("L110").Value = 2560.25
$ws.Range("M110").Value = -1280.7896
$ws.Range("N110").Value = -6650.25
$ws.Range("H128").Value = 45987
$ws.Range("J128").Value = 45987
$ws.Range("L128").Value = 45987
$ws.Range("N128").Value = -55947
$ws.Range("H132").Value = 33839.344
$ws.Range("I132").Value = 43232.168
$ws.Range("J132").Value = 5660.875
$ws.Range("K132").Value = 129696.504
$ws.Range("L132").Value = 16982.625
$ws.Range("M132").Value = -127166.504
$ws.Range("N132").Value = -22042.625
$ws.Range("H136").Value = 3996.182
$ws.Range("I136").Value = 3691.8
$ws.Range("K136").Value = 11075.4
$ws.Range("M136").Value = -8525.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2933
$ws.Range("I107").Value = 2933
$ws.Range("K107").Value = 2933
$ws.Range("M107").Value = -1013
$ws.Range("H134").Value = 1871.7693
$ws.Range("I134").Value = 1871.7693
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5615.3079
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -3080.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1670
$ws.Range("J22").Value = 2064.9
$ws.Range("L22").Value = 2064.9
$ws.Range("N22").Value = -2764.9
$ws.Range("H28").Value = 9643
$ws.Range("J28").Value = 9643
$ws.Range("L28").Value = 9643
$ws.Range("N28").Value = -10133
$ws.Range("H86").Value = 17724.967
$ws.Range("I86").Value = 30638
$ws.Range("J86").Value = 4811.933
$ws.Range("K86").Value = 30638
$ws.Range("L86").Value = 4811.933
$ws.Range("M86").Value = -29515
$ws.Range("N86").Value = -7057.933
$ws.Range("H89").Value = 17724.967
$ws.Range("I89").Value = 30638
$ws.Range("J89").Value = 4811.933
$ws.Range("K89").Value = 153190
$ws.Range("L89").Value = 24059.665
$ws.Range("M89").Value = -147574
$ws.Range("N89").Value = -35291.665
$ws.Range("H122").Value = 1999.5294
$ws.Range("I122").Value = 2003.6666
$ws.Range("J122").Value = 1989.6
$ws.Range("K122").Value = 6010.9998
$ws.Range("L122").Value = 5968.799999999999
$ws.Range("M122").Value = -3560.9998
$ws.Range("N122").Value = -10868.8
$ws.Range("H132").Value = 2001.5278
$ws.Range("I132").Value = 1813.3823
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 5440.1469
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -2910.1469
$ws.Range("N132").Value = -20660
$ws.Range("H134").Value = 102818.9
$ws.Range("I134").Value = 126023.625
$ws.Range("K134").Value = 378070.875
$ws.Range("M134").Value = -375535.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 62849.875
$ws.Range("J37").Value = 62849.875
$ws.Range("L37").Value = 188549.625
$ws.Range("N37").Value = -188773.625
$ws.Range("H136").Value = 4038.6667
$ws.Range("I136").Value = 4038.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12116.0001
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -7016.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 59999
$ws.Range("J15").Value = 59999
$ws.Range("L15").Value = 59999
$ws.Range("N15").Value = -60575
$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995
$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 43812.582
$ws.Range("I22").Value = 67327.266
$ws.Range("J22").Value = 4621.4443
$ws.Range("K22").Value = 67327.266
$ws.Range("L22").Value = 4621.4443
$ws.Range("M22").Value = -67032.266
$ws.Range("N22").Value = -5211.4443
$ws.Range("H27").Value = 43812.582
$ws.Range("I27").Value = 67327.266
$ws.Range("J27").Value = 4621.4443
$ws.Range("K27").Value = 67327.266
$ws.Range("L27").Value = 4621.4443
$ws.Range("M27").Value = -67220.266
$ws.Range("N27").Value = -4835.4443
$ws.Range("H42").Value = 32499.5
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H46").Value = 18649.334
$ws.Range("I46").Value = 45499.5
$ws.Range("J46").Value = 5224.25
$ws.Range("K46").Value = 45499.5
$ws.Range("L46").Value = 5224.25
$ws.Range("M46").Value = -45311.5
$ws.Range("N46").Value = -5600.25
$ws.Range("H49").Value = 32499.5
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H136").Value = 3647.72
$ws.Range("I136").Value = 3033.9375
$ws.Range("K136").Value = 9101.8125
$ws.Range("M136").Value = -6551.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 65100
$ws.Range("J75").Value = 65100
$ws.Range("L75").Value = 65100
$ws.Range("N75").Value = -66972
$ws.Range("H78").Value = 65100
$ws.Range("J78").Value = 65100
$ws.Range("L78").Value = 195300
$ws.Range("N78").Value = -204660
$ws.Range("H96").Value = 1975.8
$ws.Range("I96").Value = 1975.8
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1975.8
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -602.8
$ws.Range("H100").Value = 793.55554
$ws.Range("I100").Value = 737.5714
$ws.Range("J100").Value = 989.5
$ws.Range("K100").Value = 1475.1428
$ws.Range("L100").Value = 1979
$ws.Range("M100").Value = -934.1428000000001
$ws.Range("N100").Value = -3061
$ws.Range("H126").Value = 69525.97
$ws.Range("I126").Value = 85512.53999999999
$ws.Range("J126").Value = 5579.6665
$ws.Range("K126").Value = 256537.62
$ws.Range("L126").Value = 16738.9995
$ws.Range("M126").Value = -254067.62
$ws.Range("N126").Value = -21678.9995
$ws.Range("H131").Value = 43978.832
$ws.Range("J131").Value = 43978.832
$ws.Range("L131").Value = 43978.832
$ws.Range("N131").Value = -54058.832
$ws.Range("H132").Value = 66794.53
$ws.Range("I132").Value = 74407.60000000001
$ws.Range("J132").Value = 5890
$ws.Range("K132").Value = 223222.8
$ws.Range("L132").Value = 17670
$ws.Range("M132").Value = -220692.8
$ws.Range("N132").Value = -22730
$ws.Range("H136").Value = 4893.2617
$ws.Range("I136").Value = 5110.5884
$ws.Range("K136").Value = 15331.7652
$ws.Range("M136").Value = -12781.7652
$ws.Range("H140").Value = 71159.60000000001
$ws.Range("J140").Value = 71159.60000000001
$ws.Range("L140").Value = 71159.60000000001
$ws.Range("N140").Value = -81519.60000000001
